$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parsed mile posts")

# Update the vehicle/input data values
$ws.Range("D2").Value = 200000
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5

# Reset the view: scroll back to top and move selection to D3
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D3").Select()
